$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# ---------------------------------------------------------------------------
# The Kwok-Neumann (S:U) and Li-Neumann (old V:X) blocks each need a 4th,
# "average angle" column inserted right after their pair of reference
# readings, mirroring the existing Wu block pattern (F:H readings + I avg).
# This requires two new columns:
#   - a new V, inserted just left of the old V (Li-Neumann block shifts
#     from V:X to W:Y)
#   - a new Z, appended just right of the (now shifted) Y column
# ---------------------------------------------------------------------------
$ws.Columns("V").Insert()
$ws.Columns("Z").Insert()

# Column-insert copies formatting from the neighbouring column into every
# row of the new column, even rows that should stay completely empty.
# Strip those two stray carry-over cells back to "no cell at all".
$ws.Range("Z1").Clear()
$ws.Range("V10").Clear()

# New column V: average of the Kwok-Neumann readings (T:U) for each liquid
$ws.Range("V3:V8").Formula = "=AVERAGEA(T3:U3)"
$ws.Range("V9").Formula = "=_xlfn.STDEV.P(S3:S8)/SQRT(COUNT(S3:S8))"
$ws.Range("V9").ClearFormats()
$ws.Range("V9").Font.Bold = $true

# New column Z: average of the Li-Neumann readings (X:Y, formerly W:X) for
# each liquid
$ws.Range("Z3:Z8").Formula = "=AVERAGEA(X3:Y3)"
$ws.Range("Z9").Formula = "=_xlfn.STDEV.P(W3:W8)/SQRT(COUNT(W3:W8))"
$ws.Range("Z9").ClearFormats()
$ws.Range("Z9").Font.Bold = $true

# Restore the selection/active cell to reflect where the author ended up
# working (mirrors the shift of the old activeCell J13 -> X13 once two
# columns were inserted to its left).
$ws.Range("X13").Select()
